# Updated cryptos list (price/volume refresh + one rank swap) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (would otherwise be auto-parsed by Excel
# into a numeric value, losing trailing zeros / exact formatting) must be forced to stay text:
# format as Text, assign the literal string, then drop back to the default "Normal" style so
# no stray formatting is left behind on the cell.
$textCells = @{
    'D5' = '543.14'
    'D6' = '151.89'
    'D8' = '0.572'
    'D10' = '0.114'
    'D11' = '6.14'
    'D12' = '0.370'
    'D16' = '23.75'
    'D19' = '5.17'
    'D20' = '383.70'
    'D21' = '12.03'
    'D22' = '6.70'
    'D24' = '65.57'
    'D27' = '0.187'
    'D28' = '0.997'
    'D30' = '8.36'
    'D31' = '1.00'
    'D32' = '1.73'
    'D33' = '20.51'
    'D34' = '159.71'
    'D35' = '4.66'
    'D36' = '5.98'
    'D37' = '1.07'
    'D39' = '1.56'
    'D40' = '3.94'
    'D42' = '37.33'
    'D43' = '22.21'
    'D44' = '0.665'
    'D45' = '0.0594'
    'D46' = '0.0248'
    'D47' = '0.997'
    'D48' = '5.01'
    'D49' = '0.0956'
    'D50' = '19.88'
    'D51' = '269.20'
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}

# Remaining changed cells (text-like values Excel will not reinterpret as numbers) can be
# assigned directly.
$plainCells = @{
    'D2' = '61.621.39'
    'E2' = '  -3.72%  '
    'D3' = '2.969.42'
    'E3' = '  -5.48%  '
    'E4' = '  -0.05%  '
    'E5' = '  -4.14%  '
    'E6' = '  -6.03%  '
    'E7' = '  -0.23%  '
    'E8' = '  +0.18%  '
    'D9' = '2.979.01'
    'E9' = '  -5.20%  '
    'E10' = '  -1.73%  '
    'E11' = '  -6.34%  '
    'E12' = '  -2.33%  '
    'D13' = '3.488.35'
    'E13' = '  -5.62%  '
    'E14' = '  -2.74%  '
    'D15' = '61.700.04'
    'E15' = '  -3.78%  '
    'E16' = '  -4.45%  '
    'D17' = '2.973.95'
    'E17' = '  -5.60%  '
    'E18' = '  -4.15%  '
    'E19' = '  -0.41%  '
    'E20' = '  -4.03%  '
    'E21' = '  -4.13%  '
    'E22' = '  -5.32%  '
    'E23' = '  +0.27%  '
    'E24' = '  -3.81%  '
    'E25' = '  -2.06%  '
    'D26' = '3.085.84'
    'E26' = '  -6.38%  '
    'E27' = '  -2.76%  '
    'E28' = '  -0.32%  '
    'D29' = '0.0₃0943'
    'E29' = '  -5.93%  '
    'E30' = '  -4.21%  '
    'E31' = '  +0.07%  '
    'E32' = '  -3.78%  '
    'E33' = '  -2.41%  '
    'E34' = '  +2.29%  '
    'E35' = '  -2.48%  '
    'E36' = '  -3.83%  '
    'E37' = '  -2.64%  '
    'E38' = '  -3.58%  '
    'E39' = '  -5.65%  '
    'E40' = '  -2.56%  '
    'D41' = '2.412.28'
    'E41' = '  -8.97%  '
    'E42' = '  -2.56%  '
    'E43' = '  -5.50%  '
    'E44' = '  -3.67%  '
    'E45' = '  -2.59%  '
    'B46' = 'VeChain'
    'C46' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E46' = '  -2.42%  '
    'B47' = 'FirstDigitalUSD'
    'C47' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'E47' = '  -0.35%  '
    'E48' = '  -7.69%  '
    'E49' = '  -1.75%  '
    'E50' = '  -5.09%  '
    'E51' = '  -5.99%  '
}
foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}
